$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 367, shifting existing rows 367-425 down to 368-426
$ws.Rows(367).Insert()

# Populate the new row with the new derived variable "Rx25" / der_c19_treatment
$ws.Cells.Item(367, 1).Value = "Rx25"
$ws.Cells.Item(367, 2).Value = "der_c19_treatment"
$ws.Cells.Item(367, 3).Value = "Treatments"
$ws.Cells.Item(367, 8).Value = "Yes: der_rem; der_toci; der_plasma; der_steroids_c19; der_monoclonals"
$ws.Cells.Item(367, 4).Value = 'Received "promising" COVID-19 treatment (remdesivir, tocilizumab, convalescent plasma, monoclonals)'
$ws.Cells.Item(367, 5).Value = "0 = No; 1 = Yes; 99 = Unknown"
$ws.Cells.Item(367, 6).Value = "NA (missing)"
$ws.Cells.Item(367, 7).Value = "No"
# Column I ("Used by Projects") is left blank for this new row

# New row wraps onto two lines like similar multi-line description rows
$ws.Rows(367).RowHeight = 31

# Update the Table1 range to include the newly added row (A1:I425 -> A1:I426)
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:I426"))

# Reflect the editor's final scroll/selection position
$excel.ActiveWindow.ScrollRow = 358
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("D367").Select()
